$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.336.72"
$ws.Range("E2").Value = "  -4.84%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.079.77"
$ws.Range("E3").Value = "  -6.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.01"
$ws.Range("E5").Value = "  -7.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.70"
$ws.Range("E6").Value = "  -8.61%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.073.55"
$ws.Range("E8").Value = "  -6.80%  "
$ws.Range("E9").Value = "  -6.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("E10").Value = "  -9.61%  "
$ws.Range("E11").Value = "  -11.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  -10.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.621.09"
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("E15").Value = "  -7.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.627.49"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.092.23"
$ws.Range("E17").Value = "  -6.23%  "
$ws.Range("E18").Value = "  -10.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.58"
$ws.Range("E19").Value = "  -8.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("E20").Value = "  -7.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "333.40"
$ws.Range("E22").Value = "  -10.71%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -6.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.79"
$ws.Range("E25").Value = "  -8.70%  "
$ws.Range("E26").Value = "  -4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0895"
$ws.Range("E28").Value = "  -12.68%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.61"
$ws.Range("E30").Value = "  -6.26%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -9.94%  "
$ws.Range("E33").Value = "  -9.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.84"
$ws.Range("E34").Value = "  -7.59%  "
$ws.Range("E35").Value = "  -4.81%  "
$ws.Range("E36").Value = "  -7.43%  "
$ws.Range("E37").Value = "  -9.39%  "
$ws.Range("E38").Value = "  -12.46%  "
$ws.Range("E39").Value = "  -6.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.91"
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0663"
$ws.Range("E41").Value = "  -8.13%  "
$ws.Range("E42").Value = "  -10.59%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.671"
$ws.Range("E44").Value = "  -10.03%  "
$ws.Range("E45").Value = "  -6.86%  "
$ws.Range("E46").Value = "  -7.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.237.84"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("E48").Value = "  -11.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.01"
$ws.Range("E49").Value = "  -5.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.96"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("E51").Value = "  -9.85%  "
